# Update the build-timestamp strings throughout the workbook to reflect
# the new release build time, per commit "Update for release mines - January 30".

$wb = $excel.ActiveWorkbook

$oldStamp = "January 30 2026 16.19.47 EST"
$newStamp = "February 02 2026 12.49.33 EST"

$oldVersion = "mines - January 30 (built on $oldStamp)"
$newVersion = "mines - January 30 (built on $newStamp)"

# --- "About" sheet ---
$about = $wb.Worksheets.Item("About")

$about.Range("A2").Value = "Version: $newVersion"

$about.Range("A6").Value = "Recommended Citation:  " + '"Global Energy Monitor, Coal mine boundaries and methane sources for Ironbark No. 1 Coal Mine, Australia, M0052, version ' + "'$newVersion'. (See the CC license for attribution requirements if sharing or adapting the data set.)"

# --- "Boundaries and methane sources" sheet ---
$data = $wb.Worksheets.Item("Boundaries and methane sources")

for ($row = 2; $row -le 11; $row++) {
    $cell = $data.Range("S$row")
    if ($cell.Value2 -eq $oldVersion) {
        $cell.Value = $newVersion
    }
}
